$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 97
$ws.Range("F6").Value = 7001
$ws.Range("F7").Value = 658
$ws.Range("F8").Value = 141
$ws.Range("F9").Value = 12558
$ws.Range("F10").Value = 12931
$ws.Range("F12").Value = 1294
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 911
$ws.Range("F17").Value = 194
$ws.Range("F18").Value = 1431
$ws.Range("F19").Value = 359
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 1037
$ws.Range("F25").Value = 505
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 257
$ws.Range("F29").Value = 2058
$ws.Range("F30").Value = 109
$ws.Range("F33").Value = 146
$ws.Range("I33").Value = "//i2.hdslb.com/bfs/openplatform/202407/CoTyr9hO1720665458801.jpeg"
$ws.Range("F34").Value = 46
$ws.Range("F36").Value = 3756
$ws.Range("F37").Value = 4430
$ws.Range("F38").Value = 276
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 1768
$ws.Range("F44").Value = 914
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 15
$ws.Range("F3").Value = 15
$ws.Range("F5").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("F23").Value = 75
$ws.Range("F26").Value = 49
$ws.Range("F27").Value = 2
$ws.Range("F29").Value = 6

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 49

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 593
$ws.Range("F3").Value = 10529
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 12558
$ws.Range("F9").Value = 12931
$ws.Range("F10").Value = 39
$ws.Range("F13").Value = 5477
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 1037
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 505
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 257
$ws.Range("F27").Value = 2058
$ws.Range("F30").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("I33").Value = "//i2.hdslb.com/bfs/openplatform/202407/CoTyr9hO1720665458801.jpeg"
$ws.Range("F34").Value = 46
$ws.Range("F35").Value = 101
$ws.Range("F36").Value = 3
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F40").Value = 276
$ws.Range("F41").Value = 132
$ws.Range("F43").Value = 1768
$ws.Range("F45").Value = 297
$ws.Range("F46").Value = 41
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 4297
